# Updated queries for C3DC first half testcases.
#
# The TabQuery / StatQuery SQL text stored in several cells joined tables
# using the generic ".id" columns (std.id / prt.id ...). The queries are
# updated to join on the fully-qualified id columns instead
# (std.study_id / prt.participant_id ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old / new LEFT JOIN block shared by every affected query.
$oldJoin = "LEFT JOIN `n" +
    "    df_participant prt ON std.id = prt.`"study.id`"`n" +
    "LEFT JOIN `n" +
    "    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`n" +
    "LEFT JOIN `n" +
    "    df_treatments trt ON prt.id = trt.`"participant.id`"`n" +
    "LEFT JOIN `n" +
    "    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`n" +
    "LEFT JOIN `n" +
    "    df_survival srv ON prt.id = srv.`"participant.id`"`n" +
    "LEFT JOIN `n" +
    "    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newJoin = "LEFT JOIN `n" +
    "    df_participant prt ON std.study_id = prt.`"study.study_id`"`n" +
    "LEFT JOIN `n" +
    "    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`n" +
    "LEFT JOIN `n" +
    "    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`n" +
    "LEFT JOIN `n" +
    "    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`n" +
    "LEFT JOIN `n" +
    "    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`n" +
    "LEFT JOIN `n" +
    "    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

# Every cell on Sheet1 that contains one of the affected SQL queries
# (StatQuery in C2, TabQuery in B2:B7).
$cellsToUpdate = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellAddr in $cellsToUpdate) {
    $range = $ws.Range($cellAddr)
    $text = $range.Value()
    if ($text -ne $null -and $text.Contains($oldJoin)) {
        $range.Value = $text.Replace($oldJoin, $newJoin)
    }
}

# The workbook was also left scrolled/selected at a different cell
# (row 6 at top, C7 selected) when it was last saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
